$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell E1 - "Toss Winner", formatted like the other header cells
$ws.Range("E1").Value = "Toss Winner"
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# E2 / E3 become present-but-empty cells (blank toss-winner predictions
# for the already existing matches)
$ws.Range("E2").Font.Bold = $false
$ws.Range("E3").Font.Bold = $false

# New row 4 - Delhi Capitals vs Lucknow Super Giants
$ws.Range("A4").Value = "24-03-2025"
$ws.Range("B4").Value = "Delhi Capitals vs Lucknow Super Giants"
$ws.Range("C4").Font.Bold = $false
$ws.Range("D4").Value = "Delhi Capitals"
$ws.Range("E4").Value = "Delhi Capitals"
